# Generate Report for Handoff
#
# Localization status moved from "In Translation" to "Ready for handoff":
#  - Overview!E2:F2, zh-cn!C2, de-de!C2  -> "Ready for handoff"
#  - Overview!G2 / de-de!H2 (Latest HO Xliff Generate Date / Latest Handoff
#    Datetime, both "2016-09-02 15:13:00") -> "2016-09-02 15:13:36"
#  - zh-cn!H2 (Latest Handoff Datetime, "2016-09-02 15:12:55") ->
#    "2016-09-02 15:13:32"
# Also widens the "Status" columns (Overview E:F, zh-cn C, de-de C) to fit
# the new, longer "Ready for handoff" text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Generate/handoff timestamps ---
$wsOverview.Range("G2").Value = "2016-09-02 15:13:36"
$wsZhCn.Range("H2").Value = "2016-09-02 15:13:32"
$wsDeDe.Range("H2").Value = "2016-09-02 15:13:36"

# --- Widen the Status columns to fit "Ready for handoff" ---
$wsOverview.Range("E1:F1").ColumnWidth = 17.2159881591797
$wsZhCn.Range("C1").ColumnWidth = 17.2159881591797
$wsDeDe.Range("C1").ColumnWidth = 17.2159881591797
